# close #206: Adds support for zero-sum influencing factors
#
# The "valores" sheet gains a new textual marker ("DI") used to flag
# cells whose underlying ratio is a zero-sum / indeterminate influencing
# factor. Where that marker is written, the paired numeric cell in the
# same row is forced to 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (id 1100049): the influencing factor in column Q is indeterminate,
# so Q4 is flagged "DI" and its paired value P4 is zeroed out.
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = "DI"

# Row 5 (id 1100056): same situation, mirrored in column P instead.
$ws.Range("P5").Value = "DI"
$ws.Range("Q5").Value = 0

# Restore the author's on-screen scroll position / selection.
$ws.Range("H1").Select()
$excel.ActiveWindow.ScrollColumn = 8
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("Q22").Select()
